$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-of date range) ---
$ws.Cells.Item(8,1).Value = "Volume 30   Number  3"
$ws.Cells.Item(9,3).Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# --- Crime-complaint table updates (week of 1/16/2023 - 1/22/2023) ---
$ws.Range("N36").Copy() | Out-Null
$ws.Range("N15").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15,14).Value = -100
$ws.Cells.Item(16,3).Value = 1
$ws.Cells.Item(16,4).Value = 1
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 12
$ws.Cells.Item(16,7).Value = 7
$ws.Cells.Item(16,8).Value = 71.428571428571
$ws.Cells.Item(16,10).Value = 6
$ws.Cells.Item(16,11).Value = 50
$ws.Cells.Item(16,12).Value = 200
$ws.Cells.Item(16,13).Value = 80
$ws.Cells.Item(16,14).Value = -76.923076923076
$ws.Cells.Item(17,3).Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17,4).Value = 4
$ws.Range("E16").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17,5).Value = -100
$ws.Cells.Item(17,7).Value = 5
$ws.Cells.Item(17,8).Value = -20
$ws.Cells.Item(17,10).Value = 5
$ws.Cells.Item(17,11).Value = -20
$ws.Cells.Item(17,12).Value = -33.333333333333
$ws.Cells.Item(17,13).Value = -42.857142857142
$ws.Cells.Item(17,14).Value = -75
$ws.Cells.Item(18,3).Value = 6
$ws.Cells.Item(18,4).Value = 6
$ws.Cells.Item(18,6).Value = 15
$ws.Cells.Item(18,7).Value = 16
$ws.Cells.Item(18,8).Value = -6.25
$ws.Cells.Item(18,9).Value = 12
$ws.Cells.Item(18,10).Value = 14
$ws.Cells.Item(18,11).Value = -14.285714285714
$ws.Cells.Item(18,12).Value = 71.428571428571
$ws.Cells.Item(18,13).Value = 71.428571428571
$ws.Cells.Item(18,14).Value = -84.210526315789
$ws.Cells.Item(19,3).Value = 9
$ws.Cells.Item(19,4).Value = 18
$ws.Cells.Item(19,5).Value = -50
$ws.Cells.Item(19,6).Value = 36
$ws.Cells.Item(19,7).Value = 53
$ws.Cells.Item(19,8).Value = -32.075471698113
$ws.Cells.Item(19,9).Value = 29
$ws.Cells.Item(19,10).Value = 44
$ws.Cells.Item(19,11).Value = -34.090909090909
$ws.Cells.Item(19,12).Value = 31.818181818181
$ws.Cells.Item(19,13).Value = -29.268292682926
$ws.Cells.Item(19,14).Value = -73.394495412844
$ws.Cells.Item(20,3).Value = 2
$ws.Cells.Item(20,4).Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(20,5).Value = "***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(20,7).Value = 3
$ws.Cells.Item(20,8).Value = 33.333333333333
$ws.Cells.Item(20,9).Value = 3
$ws.Cells.Item(20,11).Value = 200
$ws.Cells.Item(20,12).Value = 50
$ws.Cells.Item(20,14).Value = -92.105263157894
$ws.Cells.Item(21,3).Value = 18
$ws.Cells.Item(21,4).Value = 29
$ws.Cells.Item(21,5).Value = -37.931034482758
$ws.Cells.Item(21,6).Value = 71
$ws.Cells.Item(21,7).Value = 85
$ws.Cells.Item(21,8).Value = -16.470588235294
$ws.Cells.Item(21,9).Value = 57
$ws.Cells.Item(21,10).Value = 71
$ws.Cells.Item(21,11).Value = -19.718309859154
$ws.Cells.Item(21,12).Value = 42.5
$ws.Cells.Item(21,13).Value = -5
$ws.Cells.Item(21,14).Value = -79.569892473118
$ws.Cells.Item(22,6).Value = "'0"
$ws.Range("F14").Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22,8).Value = -100
$ws.Range("L16").Copy() | Out-Null
$ws.Range("L22").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22,12).Value = -100
$ws.Cells.Item(24,3).Value = 19
$ws.Cells.Item(24,4).Value = 19
$ws.Cells.Item(24,5).Value = 0
$ws.Cells.Item(24,6).Value = 64
$ws.Cells.Item(24,7).Value = 67
$ws.Cells.Item(24,8).Value = -4.477611940298
$ws.Cells.Item(24,9).Value = 52
$ws.Cells.Item(24,10).Value = 55
$ws.Cells.Item(24,11).Value = -5.454545454545
$ws.Cells.Item(24,12).Value = -23.529411764705
$ws.Cells.Item(24,13).Value = 40.540540540540
$ws.Cells.Item(25,3).Value = 5
$ws.Cells.Item(25,4).Value = 7
$ws.Cells.Item(25,5).Value = -28.571428571428
$ws.Cells.Item(25,6).Value = 17
$ws.Cells.Item(25,7).Value = 17
$ws.Cells.Item(25,9).Value = 13
$ws.Cells.Item(25,10).Value = 14
$ws.Cells.Item(25,11).Value = -7.142857142857
$ws.Cells.Item(25,12).Value = 116.666666666667
$ws.Cells.Item(25,13).Value = -13.333333333333
$ws.Cells.Item(26,3).Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(27,4).Value = 2
$ws.Cells.Item(27,7).Value = 3
$ws.Cells.Item(27,8).Value = -33.333333333333
$ws.Cells.Item(27,10).Value = 3
$ws.Cells.Item(27,11).Value = -33.333333333333
$ws.Cells.Item(27,12).Value = -60
$ws.Cells.Item(30,7).Value = "'0"
$ws.Range("G14").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(30,8).Value = "***.*"
$ws.Range("H14").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null
